$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(201, 9, 30, 15, 45, 30),
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(101, 9, 30, 15, 60, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(601, 9, 60, 67, 60, 42),
    @(901, 16, 15, 45, 60, 60),
    @(301, 6, 45, 30, 60, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(902, 1, 0, 0, 0, 0),
    @(701, 3, 90, 45, 97, 15),
    @(801, 3, 67, 65, 52, 45),
    @(501, 9, 52, 30, 75, 45),
    @(401, 9, 48, 67, 75, 45),
    @(1, 0, 2, 2, 2, 2),
    @(2, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(802, 0, 4, 5, 4, 0),
    @(502, 0, 4, 0, 0, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(602, 0, 0, 4, 0, 9),
    @(402, 0, 0, 4, 0, 0),
    @(702, 0, 0, 0, 4, 0),
    @(1002, 0, 0, 0, 0, 9)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}
